$d = $word.ActiveDocument

# --- Locate the target paragraphs by their (trimmed) text content so the
#     script is resilient to any paragraph-index drift. ---
$licenseInfoPara = $null
$bigLicensePara   = $null
$pdfVersionPara   = $null
$listPara         = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "License Information") {
        $licenseInfoPara = $i
    } elseif ($t -eq "This PDF version is provided under the same license.") {
        $pdfVersionPara = $i
    } elseif ($t.Contains("is based on") -and $t.Contains("Biblica Bible Dictionary")) {
        $bigLicensePara = $i
    } elseif ($t.StartsWith("Мадианитяне, Македония")) {
        $listPara = $i
    }
}

# ---------------------------------------------------------------------
# 1) Rewrite the "Ключевые термины (Biblica) (Russian) is based on: ..."
#    paragraph into the new Biblica Study Notes credit paragraph.
# ---------------------------------------------------------------------
$oldBoldText = "Ключевые термины (Biblica)"
$newBoldText = "Biblica Study Notes (Key Terms)"

$p = $d.Paragraphs($bigLicensePara)
$pRange = $p.Range
# Replace just the bold run's text (scoped Find so only this paragraph's
# occurrence is touched, leaving the other two "Ключевые термины (Biblica)"
# occurrences elsewhere in the document untouched).
$pRange.Find.Execute($oldBoldText, $true, $false, $false, $false, $false, $true, 1, $false, $newBoldText, 2) | Out-Null

# Remove everything else in the paragraph that follows the (now renamed)
# bold run, up to (but excluding) the paragraph mark.
$p = $d.Paragraphs($bigLicensePara)
$full = $p.Range
$tailStart = $full.Start + $newBoldText.Length
$tailEnd = $full.End - 1
if ($tailEnd -gt $tailStart) {
    $tailRange = $d.Range($tailStart, $tailEnd)
    $tailRange.Delete()
}

# Insert the new (non-bold) trailing text just before the paragraph mark,
# then explicitly clear Bold on that new span so it forms its own run
# distinct from the "Biblica Study Notes (Key Terms)" heading run.
$newTail = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."

$p = $d.Paragraphs($bigLicensePara)
$full = $p.Range
$insPoint = $d.Range($full.End - 1, $full.End - 1)
$insPoint.InsertBefore($newTail)

$newTailStart = $tailStart
$newTailEnd = $newTailStart + $newTail.Length
$fixRange = $d.Range($newTailStart, $newTailEnd)
$fixRange.Font.Bold = 0

# ---------------------------------------------------------------------
# 2) Delete whole paragraphs (bottom-to-top so indices already found
#    above stay valid for the ones not yet removed).
# ---------------------------------------------------------------------
$toDelete = @($listPara, $pdfVersionPara, $licenseInfoPara) | Sort-Object -Descending
foreach ($idx in $toDelete) {
    if ($idx -ne $null) {
        $d.Paragraphs($idx).Range.Delete()
    }
}
